# Updates the NATMI ligand-receptor edge table (Angpt1 -> Tie1) with the
# new TPM-based figures and adds the "Resolving-Mac" sending/target cluster.
#
# Columns A:T are, in order:
#   A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#   E Ligand-expressing cells, F Ligand detection rate,
#   G Ligand average expression value, H Ligand total expression value,
#   I Ligand derived specificity of average expression value,
#   J Ligand derived specificity of total expression value,
#   K Receptor-expressing cells, L Receptor detection rate,
#   M Receptor average expression value, N Receptor total expression value,
#   O Receptor derived specificity of average expression value,
#   P Receptor derived specificity of total expression value,
#   Q Edge average expression weight, R Edge total expression weight,
#   S Edge average expression derived specificity,
#   T Edge total expression derived specificity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Angpt1", "Tie1", "ECs", 1, 0.3333333333333333, 0.09342833333333334, 0.280285, 0.00469825789230973, 0.00469825789230973, 3, 1, 106.3675996666667, 319.102799, 0.9937770853850048, 0.9937770853850048, 9.937747557523888, 89.439728017715, 0.004669021034606659, 0.004669021034606659),
    @("ECs", "Angpt1", "Tie1", "FAPs", 1, 0.3333333333333333, 0.09342833333333334, 0.280285, 0.00469825789230973, 0.00469825789230973, 1, 0.3333333333333333, 0.015201, 0.045603, 0.0001420207424279358, 0.0001420207424279358, 0.001420204095, 0.012781836855, 0.0000006672500739837368, 0.0000006672500739837368),
    @("ECs", "Angpt1", "Tie1", "MuSCs", 1, 0.3333333333333333, 0.09342833333333334, 0.280285, 0.00469825789230973, 0.00469825789230973, 2, 0.6666666666666666, 0.5886716666666666, 1.766015, 0.005499874162639981, 0.005499874162639981, 0.05499861269722222, 0.494987514275, 0.00002583982719133366, 0.00002583982719133366),
    @("ECs", "Angpt1", "Tie1", "Resolving-Mac", 1, 0.3333333333333333, 0.09342833333333334, 0.280285, 0.00469825789230973, 0.00469825789230973, 2, 0.6666666666666666, 0.06218866666666667, 0.186566, 0.0005810197099272038, 0.0005810197099272037, 0.005810183478888889, 0.05229165131000001, 0.000002729780437752996, 0.000002729780437752995),
    @("FAPs", "Angpt1", "Tie1", "ECs", 3, 1, 14.289121, 42.867363, 0.7185612021237531, 0.7185612021237531, 3, 1, 106.3675996666667, 319.102799, 0.9937770853850048, 0.9937770853850048, 1519.89950211656, 13679.09551904904, 0.7140896571172887, 0.7140896571172887),
    @("FAPs", "Angpt1", "Tie1", "FAPs", 3, 1, 14.289121, 42.867363, 0.7185612021237531, 0.7185612021237531, 1, 0.3333333333333333, 0.015201, 0.045603, 0.0001420207424279358, 0.0001420207424279358, 0.217208928321, 1.954880354889, 0.0001020505954055254, 0.0001020505954055254),
    @("FAPs", "Angpt1", "Tie1", "MuSCs", 3, 1, 14.289121, 42.867363, 0.7185612021237531, 0.7185612021237531, 2, 0.6666666666666666, 0.5886716666666666, 1.766015, 0.005499874162639981, 0.005499874162639981, 8.411600674271666, 75.704406068445, 0.003951996189835954, 0.003951996189835954),
    @("FAPs", "Angpt1", "Tie1", "Resolving-Mac", 3, 1, 14.289121, 42.867363, 0.7185612021237531, 0.7185612021237531, 2, 0.6666666666666666, 0.06218866666666667, 0.186566, 0.0005810197099272038, 0.0005810197099272037, 0.8886213828286667, 7.997592445458, 0.0004174982212228859, 0.0004174982212228858),
    @("MuSCs", "Angpt1", "Tie1", "ECs", 3, 1, 5.503190333333333, 16.509571, 0.2767405399839373, 0.2767405399839373, 3, 1, 106.3675996666667, 319.102799, 0.9937770853850048, 0.9937770853850048, 585.3611462654699, 5268.25031638923, 0.2750184072331096, 0.2750184072331096),
    @("MuSCs", "Angpt1", "Tie1", "FAPs", 3, 1, 5.503190333333333, 16.509571, 0.2767405399839373, 0.2767405399839373, 1, 0.3333333333333333, 0.015201, 0.045603, 0.0001420207424279358, 0.0001420207424279358, 0.083653996257, 0.752885966313, 0.00003930289694842663, 0.00003930289694842663),
    @("MuSCs", "Angpt1", "Tie1", "MuSCs", 3, 1, 5.503190333333333, 16.509571, 0.2767405399839373, 0.2767405399839373, 2, 0.6666666666666666, 0.5886716666666666, 1.766015, 0.005499874162639981, 0.005499874162639981, 3.239572225507222, 29.156150029565, 0.001522038145612693, 0.001522038145612693),
    @("MuSCs", "Angpt1", "Tie1", "Resolving-Mac", 3, 1, 5.503190333333333, 16.509571, 0.2767405399839373, 0.2767405399839373, 2, 0.6666666666666666, 0.06218866666666667, 0.186566, 0.0005810197099272038, 0.0005810197099272037, 0.3422360692428889, 3.080124623186, 0.000160791708266565, 0.000160791708266565)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowValues = $rows[$i]
    $targetRow = $startRow + $i
    for ($col = 0; $col -lt $rowValues.Count; $col++) {
        $ws.Cells.Item($targetRow, $col + 1).Value = $rowValues[$col]
    }
}
